$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 228
$ws.Cells.Item(2, 10).Value = 300
$ws.Cells.Item(2, 12).Value = 300
$ws.Cells.Item(2, 14).Value = -526
$ws.Cells.Item(17, 8).Value = 2386169
$ws.Cells.Item(17, 10).Value = 2569636.8
$ws.Cells.Item(17, 12).Value = 7708910.399999999
$ws.Cells.Item(17, 14).Value = -7709246.399999999
$ws.Cells.Item(93, 8).Value = 29571.428
$ws.Cells.Item(93, 10).Value = 29571.428
$ws.Cells.Item(93, 12).Value = 29571.428
$ws.Cells.Item(93, 14).Value = -34563.428
$ws.Cells.Item(132, 8).Value = 2703.6287
$ws.Cells.Item(132, 9).Value = 3069
$ws.Cells.Item(132, 11).Value = 9207
$ws.Cells.Item(132, 13).Value = -6677
$ws.Cells.Item(137, 8).Value = 1800.4572
$ws.Cells.Item(137, 9).Value = 1807.6072
$ws.Cells.Item(137, 11).Value = 5422.821599999999
$ws.Cells.Item(137, 13).Value = -2872.821599999999
$ws.Cells.Item(138, 8).Value = 20835534
$ws.Cells.Item(138, 9).Value = 30304186
$ws.Cells.Item(138, 10).Value = 4498.467
$ws.Cells.Item(138, 11).Value = 90912558
$ws.Cells.Item(138, 12).Value = 13495.401
$ws.Cells.Item(138, 13).Value = -90907418
$ws.Cells.Item(138, 14).Value = -23775.401

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 288.85715
$ws.Cells.Item(5, 9).Value = 298
$ws.Cells.Item(5, 10).Value = 266
$ws.Cells.Item(5, 11).Value = 298
$ws.Cells.Item(5, 12).Value = 266
$ws.Cells.Item(5, 13).Value = -186
$ws.Cells.Item(5, 14).Value = -490
$ws.Cells.Item(32, 8).Value = 2915.5684
$ws.Cells.Item(32, 9).Value = 2546.747
$ws.Cells.Item(32, 10).Value = 6926.5
$ws.Cells.Item(32, 11).Value = 2546.747
$ws.Cells.Item(32, 12).Value = 6926.5
$ws.Cells.Item(32, 13).Value = -2259.747
$ws.Cells.Item(32, 14).Value = -7500.5
$ws.Cells.Item(74, 8).Value = 52634660
$ws.Cells.Item(74, 9).Value = 62503220
$ws.Cells.Item(74, 11).Value = 62503220
$ws.Cells.Item(74, 13).Value = -62502346
$ws.Cells.Item(77, 8).Value = 52634660
$ws.Cells.Item(77, 9).Value = 62503220
$ws.Cells.Item(77, 11).Value = 312516100
$ws.Cells.Item(77, 13).Value = -312511732
$ws.Cells.Item(132, 8).Value = 15424
$ws.Cells.Item(132, 9).Value = 1823.7931
$ws.Cells.Item(132, 10).Value = 64724.75
$ws.Cells.Item(132, 11).Value = 5471.379300000001
$ws.Cells.Item(132, 12).Value = 194174.25
$ws.Cells.Item(132, 13).Value = -2941.379300000001
$ws.Cells.Item(132, 14).Value = -199234.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 288.85715
$ws.Cells.Item(4, 9).Value = 298
$ws.Cells.Item(4, 10).Value = 266
$ws.Cells.Item(4, 11).Value = 298
$ws.Cells.Item(4, 12).Value = 266
$ws.Cells.Item(4, 13).Value = -183
$ws.Cells.Item(4, 14).Value = -496
$ws.Cells.Item(99, 8).Value = 950.0417
$ws.Cells.Item(99, 9).Value = 929.9375
$ws.Cells.Item(99, 11).Value = 929.9375
$ws.Cells.Item(99, 13).Value = 568.0625
$ws.Cells.Item(134, 8).Value = 2353.463
$ws.Cells.Item(134, 9).Value = 2369.566
$ws.Cells.Item(134, 11).Value = 7108.697999999999
$ws.Cells.Item(134, 13).Value = -4573.697999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 1900
$ws.Cells.Item(4, 10).Value = 1900
$ws.Cells.Item(4, 12).Value = 1900
$ws.Cells.Item(4, 14).Value = -2124
$ws.Cells.Item(31, 8).Value = 3940.658
$ws.Cells.Item(31, 9).Value = 2523.077
$ws.Cells.Item(31, 10).Value = 7012.0835
$ws.Cells.Item(31, 11).Value = 2523.077
$ws.Cells.Item(31, 12).Value = 7012.0835
$ws.Cells.Item(31, 13).Value = -2228.077
$ws.Cells.Item(31, 14).Value = -7602.0835
$ws.Cells.Item(34, 8).Value = 3940.658
$ws.Cells.Item(34, 9).Value = 2523.077
$ws.Cells.Item(34, 10).Value = 7012.0835
$ws.Cells.Item(34, 11).Value = 2523.077
$ws.Cells.Item(34, 12).Value = 7012.0835
$ws.Cells.Item(34, 13).Value = -2321.077
$ws.Cells.Item(34, 14).Value = -7416.0835
$ws.Cells.Item(132, 8).Value = 1982.0426
$ws.Cells.Item(132, 9).Value = 1366.909
$ws.Cells.Item(132, 11).Value = 4100.727000000001
$ws.Cells.Item(132, 13).Value = -1570.727000000001
$ws.Cells.Item(134, 8).Value = 826.9
$ws.Cells.Item(134, 9).Value = 701.69696
$ws.Cells.Item(134, 10).Value = 1417.1428
$ws.Cells.Item(134, 11).Value = 2105.09088
$ws.Cells.Item(134, 12).Value = 4251.428400000001
$ws.Cells.Item(134, 13).Value = 429.9091200000003
$ws.Cells.Item(134, 14).Value = -9321.428400000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 4477.5
$ws.Cells.Item(3, 9).Value = 2025.5555
$ws.Cells.Item(3, 10).Value = 11833.333
$ws.Cells.Item(3, 11).Value = 6076.666499999999
$ws.Cells.Item(3, 12).Value = 35499.999
$ws.Cells.Item(3, 13).Value = -5964.666499999999
$ws.Cells.Item(3, 14).Value = -35723.999
$ws.Cells.Item(109, 8).Value = 2157.2727
$ws.Cells.Item(109, 9).Value = 1216.6666
$ws.Cells.Item(109, 10).Value = 3286
$ws.Cells.Item(109, 11).Value = 3649.9998
$ws.Cells.Item(109, 12).Value = 9858
$ws.Cells.Item(109, 13).Value = -2609.9998
$ws.Cells.Item(109, 14).Value = -11938
$ws.Cells.Item(113, 8).Value = 479.4
$ws.Cells.Item(113, 10).Value = 499.25
$ws.Cells.Item(113, 12).Value = 1497.75
$ws.Cells.Item(113, 14).Value = -5837.75
$ws.Cells.Item(131, 8).Value = 167476.1
$ws.Cells.Item(131, 10).Value = 189505.78
$ws.Cells.Item(131, 12).Value = 568517.34
$ws.Cells.Item(131, 14).Value = -578597.34

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(97, 8).Value = 1022.2308
$ws.Cells.Item(97, 9).Value = 1058.0952
$ws.Cells.Item(97, 10).Value = 871.6
$ws.Cells.Item(97, 11).Value = 1058.0952
$ws.Cells.Item(97, 12).Value = 871.6
$ws.Cells.Item(97, 13).Value = -562.0952
$ws.Cells.Item(97, 14).Value = -1863.6
$ws.Cells.Item(102, 8).Value = 12822023
$ws.Cells.Item(102, 9).Value = 14287011
$ws.Cells.Item(102, 11).Value = 14287011
$ws.Cells.Item(102, 13).Value = -14285389
$ws.Cells.Item(132, 8).Value = 22587.846
$ws.Cells.Item(132, 9).Value = 3360.261
$ws.Cells.Item(132, 10).Value = 169999.33
$ws.Cells.Item(132, 11).Value = 10080.783
$ws.Cells.Item(132, 12).Value = 509997.99
$ws.Cells.Item(132, 13).Value = -7550.782999999999
$ws.Cells.Item(132, 14).Value = -515057.99

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4234.773
$ws.Cells.Item(40, 9).Value = 3739.7368
$ws.Cells.Item(40, 11).Value = 3739.7368
$ws.Cells.Item(40, 13).Value = -3603.7368
$ws.Cells.Item(55, 8).Value = 188.76923
$ws.Cells.Item(55, 9).Value = 141.61539
$ws.Cells.Item(55, 10).Value = 235.92308
$ws.Cells.Item(55, 11).Value = 141.61539
$ws.Cells.Item(55, 12).Value = 235.92308
$ws.Cells.Item(55, 13).Value = 31.38461000000001
$ws.Cells.Item(55, 14).Value = -581.92308
$ws.Cells.Item(68, 8).Value = 2340.5
$ws.Cells.Item(68, 9).Value = 2310.2
$ws.Cells.Item(68, 10).Value = 2391
$ws.Cells.Item(68, 11).Value = 2310.2
$ws.Cells.Item(68, 12).Value = 2391
$ws.Cells.Item(68, 13).Value = -1561.2
$ws.Cells.Item(68, 14).Value = -3889
$ws.Cells.Item(71, 8).Value = 2340.5
$ws.Cells.Item(71, 9).Value = 2310.2
$ws.Cells.Item(71, 10).Value = 2391
$ws.Cells.Item(71, 11).Value = 11551
$ws.Cells.Item(71, 12).Value = 11955
$ws.Cells.Item(71, 13).Value = -7807
$ws.Cells.Item(71, 14).Value = -19443
$ws.Cells.Item(82, 8).Value = 2649.8572
$ws.Cells.Item(82, 9).Value = 2546
$ws.Cells.Item(82, 11).Value = 2546
$ws.Cells.Item(82, 13).Value = -2185
$ws.Cells.Item(85, 8).Value = 2649.8572
$ws.Cells.Item(85, 9).Value = 2546
$ws.Cells.Item(85, 11).Value = 2546
$ws.Cells.Item(85, 13).Value = -1298
$ws.Cells.Item(100, 8).Value = 1720.6
$ws.Cells.Item(100, 9).Value = 860.4286
$ws.Cells.Item(100, 10).Value = 2473.25
$ws.Cells.Item(100, 11).Value = 860.4286
$ws.Cells.Item(100, 12).Value = 2473.25
$ws.Cells.Item(100, 13).Value = -319.4286
$ws.Cells.Item(100, 14).Value = -3555.25
$ws.Cells.Item(122, 8).Value = 983380
$ws.Cells.Item(122, 9).Value = 1403021.5
$ws.Cells.Item(122, 10).Value = 4216.5
$ws.Cells.Item(122, 11).Value = 4209064.5
$ws.Cells.Item(122, 12).Value = 12649.5
$ws.Cells.Item(122, 13).Value = -4206614.5
$ws.Cells.Item(122, 14).Value = -17549.5
$ws.Cells.Item(132, 8).Value = 3199.889
$ws.Cells.Item(132, 9).Value = 1969.2307
$ws.Cells.Item(132, 11).Value = 5907.6921
$ws.Cells.Item(132, 13).Value = -3377.6921
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 770.6667
$ws.Cells.Item(136, 9).Value = 783.0526
$ws.Cells.Item(136, 10).Value = 300
$ws.Cells.Item(136, 11).Value = 2349.1578
$ws.Cells.Item(136, 12).Value = 900
$ws.Cells.Item(136, 13).Value = 200.8422
$ws.Cells.Item(136, 14).Value = -6000

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 70003
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 70003
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 70003
$ws.Cells.Item(2, 13).ClearContents()
$ws.Cells.Item(2, 14).Value = -70227
$ws.Cells.Item(54, 8).Value = 13333.333
$ws.Cells.Item(54, 10).Value = 13333.333
$ws.Cells.Item(54, 12).Value = 13333.333
$ws.Cells.Item(54, 14).Value = -14373.333
$ws.Cells.Item(113, 8).Value = 2080604.5
$ws.Cells.Item(113, 9).Value = 2153.4443
$ws.Cells.Item(113, 11).Value = 6460.3329
$ws.Cells.Item(113, 13).Value = -4290.3329
$ws.Cells.Item(136, 8).Value = 19610184
$ws.Cells.Item(136, 9).Value = 29412878
$ws.Cells.Item(136, 10).Value = 4800.5884
$ws.Cells.Item(136, 11).Value = 88238634
$ws.Cells.Item(136, 12).Value = 14401.7652
$ws.Cells.Item(136, 13).Value = -88236084
$ws.Cells.Item(136, 14).Value = -19501.7652
